# "new - cash in record"
# A new deposit ("入金") of 7000 on 2015-04-12 (serial 42068) is recorded in
# the "本金" (Principal) sheet, row 6, pushing the running balance forward.
# The previously-blank placeholder rows 7-10 get explicit 0 amounts in the
# C/D/E (入金/出金/分红) columns so the running-balance formulas keep working,
# and the 本金 sheet becomes the active / selected sheet (it was 损益 before).

$wb = $excel.ActiveWorkbook

$wsProfit = $wb.Worksheets.Item("损益")
$wsPrincipal = $wb.Worksheets.Item("本金")

# --- Row 6: turn the blank placeholder row into a real cash-in entry ------
# Clone the formatting of the row above (dates/currency number formats)
# so the new cells pick up the same style already used by rows 2-5,
# instead of inventing brand-new style entries.
$wsPrincipal.Range("B5:F5").Copy() | Out-Null
$wsPrincipal.Range("B6:F6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$wsPrincipal.Range("B6").Value = 42068
$wsPrincipal.Range("C6").Value = 7000
$wsPrincipal.Range("D6").Value = 0
$wsPrincipal.Range("E6").Value = 0
$wsPrincipal.Range("F6").Formula = "=F5+C6-D6+E6"

# --- Rows 7-10: fill the 入金/出金/分红 columns with explicit zeros --------
# (amount columns only; date (B) and running-total (F) stay blank, matching
# the template layout of the still-unused rows)
$wsPrincipal.Range("C7:E7").Value = 0
$wsPrincipal.Range("C8:E8").Value = 0
$wsPrincipal.Range("C9:E9").Value = 0
$wsPrincipal.Range("C10:E10").Value = 0

# --- Make "本金" the active / selected sheet instead of "损益" -------------
$wsPrincipal.Activate() | Out-Null
$wsPrincipal.Range("G6").Select() | Out-Null
